$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace letter labels (U / L) with picture filenames (p1.jpg / p2.jpg) ---
# Rows 2-11 hold the trial data in columns A and B; whichever cell currently
# reads "U" becomes "p1.jpg", whichever reads "L" becomes "p2.jpg".
for ($r = 2; $r -le 11; $r++) {
    foreach ($col in @("A", "B")) {
        $cell = $ws.Range($col + $r)
        $cur = $cell.Value()
        if ($cur -eq "U") {
            $cell.Value = "p1.jpg"
        } elseif ($cur -eq "L") {
            $cell.Value = "p2.jpg"
        }
    }
}

# --- Re-style those same cells with the new (non-hiragana) font ---
$dataRange = $ws.Range("A2:B11")
$dataRange.Font.Name = "Calibri "
$dataRange.Font.Size = 11

# --- Selection moves from C11 to B11 ---
$ws.Range("B11").Select()

# --- Page setup: portrait, paper size 9 (A4) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "done"
